$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F: header "time_taken" matching the existing header style
$ws.Range("F1").Value = "time_taken"
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$timestamps = @(
    "2021-10-05 10:52:56.868043",
    "2021-10-05 10:52:56.868055",
    "2021-10-05 10:52:56.868059",
    "2021-10-05 10:52:56.868062",
    "2021-10-05 10:52:56.868066",
    "2021-10-05 10:52:56.868069",
    "2021-10-05 10:52:56.868072",
    "2021-10-05 10:52:56.868075",
    "2021-10-05 10:52:56.868078",
    "2021-10-05 10:52:56.868081",
    "2021-10-05 10:52:56.868084",
    "2021-10-05 10:52:56.868087",
    "2021-10-05 10:52:56.868090",
    "2021-10-05 10:52:56.868093",
    "2021-10-05 10:52:56.868096"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
